# Insert a new row for "chemical_recycling_pyrolysis" right after
# "chemical_recycling_gasification" (which lives on row 9), shifting all
# subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = "chemical_recycling_pyrolysis"
$ws.Range("B10").Value = $true
